# Insert a new price-report record as row 19 (pushing the existing rows
# 19-104 down to 20-105), matching the source data's weekly update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 10
$ws.Cells.Item(19, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value = "La Araucanía"
$ws.Cells.Item(19, 4).Value = 45250
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(19, 6).Value = 100112026
$ws.Cells.Item(19, 7).Value = "Haba"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 50
$ws.Cells.Item(19, 11).Value = 10000
$ws.Cells.Item(19, 12).Value = 10000
$ws.Cells.Item(19, 13).Value = 10000
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Región del Maule"
$ws.Cells.Item(19, 16).Value = 400
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
